$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 2 and row 14 values for columns A and C (B is identical in both rows)
$a2 = $ws.Range("A2").Value2
$c2 = $ws.Range("C2").Value2
$a14 = $ws.Range("A14").Value2
$c14 = $ws.Range("C14").Value2

$ws.Range("A2").Value2 = $a14
$ws.Range("C2").Value2 = $c14
$ws.Range("A14").Value2 = $a2
$ws.Range("C14").Value2 = $c2
